# Daily attendance processing - swap the order of the two comma-separated
# names/identities recorded in the "Recorded By" column (column G) whenever
# "dnasr281@gmail.com" is listed first, e.g.
#   "dnasr281@gmail.com, System"       -> "System, dnasr281@gmail.com"
#   "dnasr281@gmail.com, admin@admin.com" -> "admin@admin.com, dnasr281@gmail.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Value2

    if ($text -ne $null -and $text -like "dnasr281@gmail.com, *") {
        $parts = $text -split ", ", 2
        $newText = $parts[1] + ", " + $parts[0]
        $cell.Value = $newText
    }
}
